{"js": "// Update the division-problem worksheet values in place.\n// The document contains a single table; every 4th row (0, 4, 8, 12, 16)\n// holds 5 division problems (\"NN\u00f7N=\") in its cells, the other rows are\n// blank spacer rows. We replace the problems, by position, with the new\n// values from the target revision - text-matching alone is unreliable\n// because several source problems repeat (e.g. \"43\u00f73=\" appears twice with\n// different replacements).\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Ordered list of the non-blank rows' new values (row index -> 5 cell values).\nconst newValuesByRow = {\n  0: [\"47\u00f79=\", \"12\u00f78=\", \"83\u00f75=\", \"30\u00f77=\", \"55\u00f76=\"],\n  4: [\"71\u00f75=\", \"81\u00f74=\", \"17\u00f78=\", \"40\u00f79=\", \"72\u00f74=\"],\n  8: [\"28\u00f77=\", \"87\u00f79=\", \"61\u00f78=\", \"59\u00f78=\", \"14\u00f72=\"],\n  12: [\"45\u00f73=\", \"31\u00f78=\", \"13\u00f72=\", \"47\u00f74=\", \"16\u00f78=\"],\n  16: [\"21\u00f72=\", \"28\u00f78=\", \"47\u00f74=\", \"14\u00f75=\", \"97\u00f75=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const row = rows.items[rowIndex];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  const values = newValuesByRow[rowIndex];\n  for (let c = 0; c < cells.items.length; c++) {\n    cells.items[c].value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet values in place.\n# The document contains a single table; every 4th row (1, 5, 9, 13, 17 in\n# Word's 1-based indexing) holds 5 division problems (\"NN\u00f7N=\") in its\n# cells, the other rows are blank spacer rows. We replace the problems, by\n# position, with the new values from the target revision - text-matching\n# alone is unreliable because several source problems repeat (e.g.\n# \"43\u00f73=\" appears twice with different replacements).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValuesByRow = @{\n    1  = @(\"47\u00f79=\", \"12\u00f78=\", \"83\u00f75=\", \"30\u00f77=\", \"55\u00f76=\")\n    5  = @(\"71\u00f75=\", \"81\u00f74=\", \"17\u00f78=\", \"40\u00f79=\", \"72\u00f74=\")\n    9  = @(\"28\u00f77=\", \"87\u00f79=\", \"61\u00f78=\", \"59\u00f78=\", \"14\u00f72=\")\n    13 = @(\"45\u00f73=\", \"31\u00f78=\", \"13\u00f72=\", \"47\u00f74=\", \"16\u00f78=\")\n    17 = @(\"21\u00f72=\", \"28\u00f78=\", \"47\u00f74=\", \"14\u00f75=\", \"97\u00f75=\")\n}\n\nforeach ($rowIndex in $newValuesByRow.Keys) {\n    $values = $newValuesByRow[$rowIndex]\n    for ($c = 1; $c -le $values.Count; $c++) {\n        $cell = $t.Cell($rowIndex, $c)\n        $cell.Range.Text = $values[$c - 1]\n    }\n}\n"}
